$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "乙肝表面抗原"
$ws.Range("A3").Value = "乙肝表面抗原"
$ws.Range("A4").Value = "乙肝表面抗原"
$ws.Range("A5").Value = "乙肝e抗体"
$ws.Range("A6").Value = "乙型肝炎核心抗体定量"
$ws.Range("A7").Value = "乙肝病毒e抗原"
$ws.Range("B7").Value = "（酶免法）0.305"
$ws.Range("C7").Value = "<1.000"
